$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 202
$ws.Range("J2").Value = 882
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 226
$ws.Range("N2").Value = 164
$ws.Range("S2").Value = 101
$ws.Range("T2").Value = 138
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 1344
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 1395
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 25
$ws.Range("AA2").Value = 2
